# Update to support full month names and Resources column format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Group" column (E) is no longer needed — remove it entirely so every
# column to its right (Data, Location, month columns, ...) shifts one spot
# to the left, which also drops the now-unused last column (old "Dec"/S).
$ws.Range("E1:E11").EntireColumn.Delete()

# Header tweaks: "Task1" -> "Task 1" and "Resource" -> "Resources".
$ws.Range("B1").Value = "Task 1"
$ws.Range("D1").Value = "Resources"

# Expand the abbreviated month headers (now sitting in G1:R1 after the
# column shift) to their full names.
$months = @("January","February","March","April","May","June","July","August","September","October","November","December")
for ($i = 0; $i -lt $months.Length; $i++) {
    $col = 7 + $i  # G is column 7
    $ws.Cells.Item(1, $col).Value = $months[$i]
}
